$wb = $excel.ActiveWorkbook

# --- Rename "Include from v3-PatientImport" sheet to "Include ValueSets" ---
$wsInclude = $wb.Worksheets.Item(2)
$wsInclude.Name = "Include ValueSets"

# --- Update the Metadata sheet ---
$wsMeta = $wb.Worksheets.Item(1)

# Version: 5.0.0 -> 6.0.0
$wsMeta.Range("B3").Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$wsMeta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value: (empty) -> Alvearie Team
$wsMeta.Range("B9").Value = "Alvearie Team"

# Row 10 becomes Jurisdiction / United States of America (was Contact / No display for ContactDetail)
$wsMeta.Range("A10").Value = "Jurisdiction"
$wsMeta.Range("B10").Value = "United States of America"

# Row 11 becomes Description / Examples of customer-specific patient status codes
$wsMeta.Range("A11").Value = "Description"
$wsMeta.Range("B11").Value = "Examples of customer-specific patient status codes"

# Row 12 becomes Purpose / (empty)
$wsMeta.Range("A12").Value = "Purpose"
$wsMeta.Range("B12").Value = ""

# Row 13 becomes Copyright / (empty)
$wsMeta.Range("A13").Value = "Copyright"
$wsMeta.Range("B13").Value = ""

# Row 14 becomes Immutable / BooleanType[null]
$wsMeta.Range("A14").Value = "Immutable"
$wsMeta.Range("B14").Value = "BooleanType[null]"

# Remove the old row 15 (previously held Immutable / BooleanType[null])
$wsMeta.Rows.Item(15).Delete()

# --- Update the "Include ValueSets" sheet (formerly Include from v3-PatientImport) ---
# Replace the code-table rows with a simple ValueSet URL reference
$wsInclude.Range("A1").Value = "ValueSet URL"
$wsInclude.Range("A2").Value = "http://terminology.hl7.org/ValueSet/v3-PatientImportance"
$wsInclude.Rows.Item(4).Delete()
$wsInclude.Rows.Item(3).Delete()
